$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new shared-string text blocks (kept as here-strings so the
#     Cypher backtick-quoted column aliases don't need escaping) ---
$samplesTab = @'
SamplesTab
'@
$filesTab = @'
FilesTab
'@
$samplesQuery = @'
MATCH (ss:study_subject)
WITH COLLECT(ss.study_subject_id) AS all_subjects
MATCH (samp:sample)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (samp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
WHERE s.study_acronym IN ["C"]  
WITH
    distinct lp,
    toInteger(split(ss.study_subject_id,'-')[2]) AS subject_id_num,
    collect(distinct f.file_id) AS files,
    samp, ss, s, p, all_subjects
RETURN
 samp.sample_id AS `Sample ID`,
            ss.study_subject_id AS `Case ID`,
            p.program_acronym AS `Program Code`,
            s.study_acronym AS `Arm`,
            ss.disease_subtype AS `Diagnosis`,
            samp.tissue_type AS `Tissue Type`,
            samp.composition AS `Tissue Composition`,
            samp.sample_anatomic_site AS `Sample Anatomic Site`,
            samp.method_of_sample_procurement AS `Sample Procurement Method`
'@
$filesQuery = @'
MATCH (f:file)-->(parent)
MATCH (f)-[:file_of_sample]->(samp)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)-[:study_subject_of_study]->(s)
MATCH (s)-[:study_of_program]->(p)
MATCH (d)-[:diagnosis_of_study_subject]->(ss)
MATCH (tp)-[:tp_of_diagnosis]->(d)
WHERE s.study_acronym IN ["C"]  
WITH
        f, parent,p, ss, d,tp, s, samp,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, parent,p, ss, d,tp, s, samp,
        f.file_size /(1024^i) AS value,
        10^precision AS factor,
        units[i] as unit
WITH
        f, parent,p, ss, d,tp, s, samp, unit,
        round(factor * value)/factor AS size
RETURN Distinct
    f.file_name AS `File Name`,
    head(labels(samp)) AS `Association`,
    f.file_description AS `Description`,
    f.file_format AS `File Format`,
     CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
    p.program_acronym AS `Program Code`,
    s.study_acronym AS `Arm`,
    ss.study_subject_id AS `Case ID`,
    samp.sample_id AS `Sample ID`
    order by f.file_name
'@
$statQuery = @'
MATCH (ss:study_subject)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)

WITH ss
MATCH (ss)-[:study_subject_of_study]->(s)
MATCH (s)-[:study_of_program]->(p)
MATCH (ss)<-[:sample_of_study_subject]-(samp)
MATCH (samp)<-[:file_of_sample]-(f)
MATCH (lp)<-[:file_of_laboratory_procedure]-(f)
WHERE s.study_acronym IN ["C"]  
RETURN COUNT(DISTINCT p) AS Programs,
COUNT(DISTINCT s) AS Arms,
COUNT(DISTINCT ss) AS Cases,
COUNT(DISTINCT samp) AS Samples,
COUNT(DISTINCT lp) AS Assays,
COUNT(DISTINCT f) AS Files
'@
$neo4jFile = @'
TC03_Bento_Filter_Arm-C_Neo4jData.xlsx
'@
$webFile = @'
TC03_Bento_Filter_Arm-C_WebData.xlsx
'@

# --- row 3 : SamplesTab ---
# Column A and column B are written first (in that order across both new
# rows) so the new unique shared-string entries land at the same indices
# the recorded workbook uses (10=SamplesTab, 11=FilesTab, 12=samples
# query, 13=files query); C/D/E below just reuse existing shared strings.
$ws.Range("A3").Value = $samplesTab
$ws.Range("A4").Value = $filesTab

$ws.Range("B3").Value = $samplesQuery
$ws.Range("B3").WrapText = $true
$ws.Range("B4").Value = $filesQuery
$ws.Range("B4").WrapText = $true

$ws.Range("C3").Value = $statQuery
$ws.Range("C3").WrapText = $true
$ws.Range("D3").Value = $neo4jFile
$ws.Range("E3").Value = $webFile

# --- row 4 : FilesTab ---
$ws.Range("C4").Value = $statQuery
$ws.Range("C4").WrapText = $true
$ws.Range("D4").Value = $neo4jFile
$ws.Range("E4").Value = $webFile

# --- row heights for the two new rows ---
$ws.Rows.Item(3).RowHeight = 360
$ws.Rows.Item(4).RowHeight = 409.6

# --- row 2 height shrank slightly in the saved workbook ---
$ws.Rows.Item(2).RowHeight = 316.8

# --- column widths were re-autofit for the new, wider content ---
$ws.Columns.Item(1).ColumnWidth = 12.666666666666666
$ws.Columns.Item(2).ColumnWidth = 75
$ws.Columns.Item(3).ColumnWidth = 50.666666666666664
$ws.Columns.Item(4).ColumnWidth = 44
$ws.Columns.Item(5).ColumnWidth = 42

# --- view state: zoom in, and leave C2:E4 selected ---
$ws.Application.ActiveWindow.Zoom = 60
$ws.Range("C2:E4").Select() | Out-Null
